$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Column headers: roll the reporting periods forward by one quarter ---
$periods = @(
"فصل سوم منتهی به 1399/09",  "فصل چهارم منتهی به 1399/12",  "فصل اول منتهی به 1400/03",  "فصل دوم منتهی به 1400/06",  "فصل سوم منتهی به 1400/09",  "فصل چهارم منتهی به 1400/12",  "فصل اول منتهی به 1401/03",  "فصل دوم منتهی به 1401/06",  "فصل سوم منتهی به 1401/09",  "فصل چهارم منتهی به 1401/12"
)
$dates = @(
"1399-10-27",  "1401-04-01 (8)",  "1400-04-16",  "1400-08-25 (2)",  "1400-10-29",  "1402-02-30 (8)",  "1401-04-28",  "1401-08-25 (2)",  "1401-10-29",  "1402-02-30 (2)"
)
$cols = @("D","E","F","G","H","I","J","K","L","M")
for ($i = 0; $i -lt 10; $i++) {
    $ws.Range($cols[$i] + "8").Value = $periods[$i]
    $ws.Range($cols[$i] + "9").Value = $dates[$i]
}

# --- Balance-sheet figures: shift every quarterly column left one slot and append the new quarter ---
$rowData = @{
    12 = @(69417, 50630, 69044, 27220, 199788, 124711, 69166, 548992, 41389, 255126)
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    14 = @(715852, 768470, 923574, 1322766, 1417184, 1298799, 2115298, 3035354, 4658376, 4133142)
    15 = @(350186, 546798, 567628, 644099, 685158, 680684, 1389433, 1206243, 1373110, 1509807)
    16 = @(162688, 71616, 153655, 147025, 204791, 693571, 209930, 105015, 312655, 468936)
    17 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    18 = @(1298143, 1437514, 1713901, 2141110, 2506921, 2797765, 3783827, 4895604, 6385530, 6367011)
    19 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    20 = @(69, 69, 69, 69, 69, 69, 69, 69, 69, 69)
    21 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    22 = @(4498746, 4473707, 4452375, 4420163, 4382650, 4485773, 4448774, 4441313, 4431833, 4440095)
    23 = @(1302, 1302, 1302, 1302, 1302, 1302, 1302, 1302, 1302, 1302)
    24 = @(40, 40, 40, 40, 40, 40, 40, 40, 40, 40)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(4500117, 4475078, 4453746, 4421534, 4384021, 4487144, 4450145, 4442684, 4433204, 4441466)
    27 = @(5798260, 5912592, 6167647, 6562644, 6890942, 7284909, 8233972, 9338288, 10818734, 10808477)
    29 = @(276627, 214295, 286806, 422459, 529206, 327294, 773044, 812952, 878852, 803973)
    30 = @(40, 40, 40, 40, 40, 40, 40, 40, 40, 40)
    31 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    32 = @(84420, 114998, 129858, 147534, 172675, 147625, 227816, 364392, 472171, 416996)
    33 = @(57745, 5323, 5288, 229475, 133482, 8819, 8724, 455208, 263834, 11132)
    34 = @(41498, 42755, 44051, 45386, 38575, 532549, 402903, 416876, 1078754, 683417)
    35 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    36 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    37 = @(460290, 377371, 466003, 844854, 873938, 1016287, 1412487, 2049428, 2693611, 1915518)
    38 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    39 = @(40, 40, 40, 40, 40, 40, 40, 40, 40, 40)
    40 = @(38575, 27403, 15892, 4032, 0, 0, 0, 0, 0, 0)
    41 = @(129351, 130856, 147507, 170541, 221949, 250921, 291630, 332557, 365179, 406132)
    42 = @(167926, 158259, 163399, 174573, 221949, 250921, 291630, 332557, 365179, 406132)
    43 = @(628216, 535630, 629402, 1019427, 1095887, 1267208, 1704117, 2381985, 3058790, 2321650)
    45 = @(4484000, 4484000, 4484000, 4484000, 4484000, 4484000, 4484000, 4484000, 4484000, 4484000)
    46 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    47 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    48 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    49 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    50 = @(43660, 54006, 62070, 73529, 86121, 97253, 122861, 166603, 206785, 243129)
    51 = @(12895, 12895, 12895, 12895, 12895, 12895, 12895, 12895, 12895, 12895)
    52 = @(40, 40, 40, 40, 40, 40, 40, 40, 40, 40)
    53 = @(18386, 18350, 18187, 18025, 17862, 17827, 17696, 17502, 17042, 17304)
    54 = @(40, 40, 40, 40, 40, 40, 40, 40, 40, 40)
    55 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    56 = @(611103, 807711, 961093, 954768, 1194177, 1405726, 1892403, 2275303, 3039222, 3729499)
    57 = @(5170044, 5376962, 5538245, 5543217, 5795055, 6017701, 6529855, 6956303, 7759944, 8486827)
    58 = @(5798260, 5912592, 6167647, 6562644, 6890942, 7284909, 8233972, 9338288, 10818734, 10808477)
}
foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# --- Column widths: the date/period columns shift pattern by one as well ---
$widths = @(29, 31, 29, 29, 29, 31, 29, 29, 29, 31)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Columns.Item(4 + $i).ColumnWidth = $widths[$i] - 0.83
}

# --- Row heights (title/header rows reflow slightly) ---
$ws.Rows.Item(2).RowHeight = 15.6
$ws.Rows.Item(5).RowHeight = 40.8
$ws.Rows.Item(6).RowHeight = 40.8
$ws.Rows.Item(8).RowHeight = 15.6

"edit applied"